# Update the "Secondary" education row (row 10) with refreshed percentages.
# Values are stored as text in the sheet (e.g. "0.26"), so a leading
# apostrophe is used to force text entry and avoid Excel re-interpreting
# the numeric-looking strings as numbers. The Style reset afterwards keeps
# the cell's formatting identical to before the edit (no quote-prefix
# style lingering on the cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B10" = "0.29"
    "C10" = "0.49"
    "D10" = "0.39"
    "E10" = "0.5"
    "F10" = "0.42"
    "G10" = "0.57"
    "H10" = "0.47"
    "I10" = "0.62"
    "J10" = "0.65"
    "K10" = "0.62"
    "L10" = "0.6"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
